$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Worksheet "NFTRTickets": make room for 3 new columns (U,V,W) by
#    shifting the existing U:Y block three columns to the right (->X:AB)
#    WITHOUT touching column-width metadata (<cols>), so we copy/paste
#    the block instead of doing a real Columns.Insert (which would also
#    shift the <cols> width definitions - not what the target file does).
# ------------------------------------------------------------------
$wsT = $wb.Worksheets.Item("NFTRTickets")

$wsT.Range("U1:Y5").Copy()
$wsT.Range("X1").PasteSpecial(-4104)   # xlPasteAll
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2) Populate the freshly-opened U:W columns with the new "Issue Field 5"
#    / mandatory-flag data.
# ------------------------------------------------------------------

# Row 1 - headers
$wsT.Range("U1").Value = "Issue Field 5 - Label"
$wsT.Range("U1").Style = $wsT.Range("T1").Style
$wsT.Range("V1").Value = "Field 5 Type"
$wsT.Range("V1").Style = $wsT.Range("T1").Style
$wsT.Range("W1").Value = "Is Mandatory "
$wsT.Range("W1").Style = $wsT.Range("T1").Style

# Row 2 - first data row gets real values
$wsT.Range("U2").Value = "Blank SIM ICCID"
$wsT.Range("U2").Style = $wsT.Range("T2").Style
$wsT.Range("V2").Value = "Text Box"
$wsT.Range("V2").Style = $wsT.Range("T2").Style
$wsT.Range("W2").Value = "Yes"
$wsT.Range("W2").Style = $wsT.Range("T2").Style

# Rows 3-5 stay blank (just inherit the bordered "empty" style like F:T)
$wsT.Range("U3:W5").ClearContents()
$wsT.Range("U3:W5").Style = $wsT.Range("T3").Style

# ------------------------------------------------------------------
# 3) View-state: NFTRTickets becomes the active tab / sheet, with an
#    updated selection. (LoginCredentials loses tabSelected, and
#    activeTab on the workbook moves to index 1 automatically.)
# ------------------------------------------------------------------
$wsT.Activate()
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
$wsT.Range("V7").Select()
